# fix the totalFee equals to NaN to newly created billing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the invoice/creation timestamp (H9)
$ws.Range("H9").Value2 = 45261.84317129629

# Clear the stray reference number that was being shown in G11 (was causing the NaN total fee)
$ws.Range("G11").Value2 = ""

# Update the billing line timestamp (B16)
$ws.Range("B16").Value2 = 45261.84322916667

# Replace the description on the first billing line with the shipment number text
$ws.Range("D16").Value2 = "SHIPMENT NO.: 123123234124"

# Fix the fee amount for the newly created billing line
$ws.Range("H16").Value2 = 143

# Keep the SPO number text on the following line
$ws.Range("D17").Value2 = "SPO NO.: 1241421414"

# Restore row heights that Excel auto-adjusted as a side effect of the value edits above,
# keeping the layout identical to the original template
$ws.Rows(11).RowHeight = 6
$ws.Rows(16).RowHeight = 9.75
